$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Pitch-by-pitch table (Pitch / Choice / Result columns F/G/H) ---
# First at-bat block (rows 10-13)
$ws.Range("F10").Value = "SL"
$ws.Range("G10").Value = "Take"
$ws.Range("H10").Value = "Ball"

$ws.Range("F11").Value = "SL"
$ws.Range("G11").Value = "Take"
$ws.Range("H11").Value = "Ball"

$ws.Range("F12").Value = "CH"
$ws.Range("G12").Value = "Take"
$ws.Range("H12").Value = "Ball"

$ws.Range("F13").Value = "FB"
$ws.Range("G13").Value = "Take"
$ws.Range("H13").Value = "Ball"

# Clear the (unused) Launch Angle value for this at-bat
$ws.Range("M12").ClearContents()

# Result of the first at-bat
$ws.Range("M15").Value = "Walk"

# Pitch mix order updated
$ws.Range("J17").Value = "CH,FB,SL"

# Second at-bat block (rows 19-22)
$ws.Range("F19").Value = "CH"
$ws.Range("G19").Value = "Take"
$ws.Range("H19").Value = "Strike"

$ws.Range("F20").Value = "CH"
$ws.Range("G20").Value = "Swing"
$ws.Range("H20").Value = "In Play"

# Exit velo / launch angle results for the ball put in play
$ws.Range("M19").Value = "94.99 MPH"
$ws.Range("M21").Value = "4.16°"

# Pitch mix order updated
$ws.Range("J26").Value = "CH,FB,SL"
